# Update gh-pages to output generated at 456a3b4
# Update "想去人数" (number of people wanting to go) figures in the
# 展览 (Exhibitions) sheet and the 全部类型 (All Types) sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - row 3 (F3) and row 4 (F4)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 238
$wsExhibit.Range("F4").Value = 861

# Sheet "全部类型" - row 4 (F4) and row 5 (F5)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 238
$wsAll.Range("F5").Value = 861
